# "LOD final working version" - add the final SVM/Vol_Name row (NetApp)
# to the LOD worksheet and leave the selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 8): SVM_Name repeats "svm1_cluster2" (row 7's
# value, reused as a shared string) paired with the new Vol_Name "NetApp".
$ws.Range("A8").Value = "svm1_cluster2"
$ws.Range("B8").Value = "NetApp"

# Move/restore the active selection to B12, matching the author's final
# cursor position when the workbook was saved.
$ws.Range("B12").Select() | Out-Null
